# Atualização de layout e nível de acesso
#
# Adds a new "AccessLevel" column (C) with a per-user access level, turns the
# placeholder row 7 (numeric 1/1) into a real "Alison" username/password
# entry, and replaces Guilherme's numeric password with the text password
# "18051980.ga". Finishes by resizing the new column and refreshing the
# sheet's active-cell selection, matching where the editor was left after
# making the changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "AccessLevel" column ------------------------------------------------
$ws.Range("C3").Value = "user"
$ws.Range("C2").Value = "coord"
$ws.Range("C4").Value = "admin"

# --- Row 7: was a placeholder (1 / 1) numeric row, now a real account ------
$ws.Range("A7").Value = "Alison"
$ws.Range("B7").Value = "Alison"
$ws.Range("B7").HorizontalAlignment = -4108  # xlCenter, matches the other password cells

$ws.Range("C1").Value = "AccessLevel"
$ws.Range("C7").Value = "diretor"

# --- Guilherme's password becomes text instead of a bare number ------------
$ws.Range("B6").Value = "18051980.ga"

# --- Remaining AccessLevel cells --------------------------------------------
$ws.Range("C5").Value = "user"
$ws.Range("C6").Value = "user"

# Size the new column to fit its contents (matches the original author's
# "best fit" column width for the AccessLevel header/values).
$ws.Columns.Item(3).ColumnWidth = 9.83

# Reflect the active selection left behind after the edit.
$ws.Range("B6").Select() | Out-Null
